# "Updated invalid inputs sheet"
# - On the "TC3 Invalid Inputs" sheet, rename the junk PartNumber test value
#   from "AR-JUNK" to "AR-JUNK2" (cell C3).
# - Move the sheet's saved cursor/selection from D7 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC3 Invalid Inputs")

$ws.Activate()

# Update the invalid PartNumber test value.
$ws.Range("C3").Value = "AR-JUNK2"

# Update the active selection saved with the sheet view.
$ws.Range("C5").Select() | Out-Null
